# Delete rows 49, 47 and 31 (in descending order so row indices of the
# rows still to be removed remain stable while earlier deletions happen).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(49).Delete()
$ws.Rows.Item(47).Delete()
$ws.Rows.Item(31).Delete()
